# Prepend a new "Jornada 33" round of fixtures above the existing schedule
# on Sheet1, pushing the current rows 1-50 ("Jornada 34"-"Jornada 38") down
# to rows 11-60. The four teams that are the "home"/"away" side of one of
# these new rows need a (new) shared string distinct from the existing
# concatenated-name strings already used elsewhere in the sheet (e.g.
# "Ath Bilbao" vs. the pre-existing "AthBilbao").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Shift everything down by 10 rows to make room for the new fixtures.
$ws.Rows("1:10").Insert()

$newFixtures = @(
    @(33, "Vallecano",   "Osasuna"),
    @(33, "Ath Bilbao",  "Barcelona"),
    @(33, "Ath Madrid",  "Real Madrid"),
    @(33, "Levante",     "Celta"),
    @(33, "Zaragoza",    "Mallorca"),
    @(33, "Espanol",     "Granada"),
    @(33, "Malaga",      "Getafe"),
    @(33, "Sociedad",    "Valencia"),
    @(33, "Valladolid",  "Sevilla"),
    @(33, "Betis",       "La Coruna")
)

for ($i = 0; $i -lt $newFixtures.Count; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $newFixtures[$i][0]
    $ws.Cells.Item($r, 2).Value = $newFixtures[$i][1]
    $ws.Cells.Item($r, 3).Value = $newFixtures[$i][2]
}

# Match the author's final selection/scroll state in the saved file.
$ws.Range("H9").Select()
